{"js": "// The reference diff for this document touches word/document.xml,\n// word/footer1-3.xml, word/footnotes.xml, word/header1-3.xml and\n// word/styles.xml, but every single changed line is a pure\n// attribute/namespace *reordering* (OOXML serializers alphabetizing\n// xmlns:* declarations and w:* attributes, e.g.\n// `w:type=\"even\" r:id=\"rId6\"` -> `r:id=\"rId6\" w:type=\"even\"`,\n// `w:top=\"1417\" w:right=\"1417\" ...` -> `w:bottom=\"1417\" w:footer=\"708\" ...`,\n// `w:styleId=\"Normal\" w:type=\"paragraph\"` -> `w:default=\"1\" w:styleId=\"Normal\" w:type=\"paragraph\"`,\n// etc.). No text, run/paragraph formatting, style definition, header/\n// footer content, footnote content, section/page setup value, or any\n// other semantically observable property differs between the before\n// and after XML - every w:val / w:id / r:id / dimension stays exactly\n// the same, only the left-to-right order of attributes/xmlns\n// declarations inside tags changes.\n//\n// That kind of change is an artifact of whichever tool re-serialized\n// the .xml parts (e.g. the document being re-saved), not something\n// exposed by the Word JavaScript API: Office.js has no attribute-order\n// / raw-OOXML-serialization control, and touching a part through the\n// object model to force a re-write would risk introducing real\n// (unwanted) content differences instead of reproducing this purely\n// cosmetic reordering. So there is nothing for this script to change -\n// we just confirm the body is reachable and leave the document exactly\n// as authored.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The reference diff for this document touches word/document.xml,\n# word/footer1-3.xml, word/footnotes.xml, word/header1-3.xml and\n# word/styles.xml, but every single changed line is a pure\n# attribute/namespace *reordering* (OOXML serializers alphabetizing\n# xmlns:* declarations and w:* attributes, e.g.\n# `w:type=\"even\" r:id=\"rId6\"` -> `r:id=\"rId6\" w:type=\"even\"`,\n# `w:top=\"1417\" w:right=\"1417\" ...` -> `w:bottom=\"1417\" w:footer=\"708\" ...`,\n# `w:styleId=\"Normal\" w:type=\"paragraph\"` -> `w:default=\"1\" w:styleId=\"Normal\" w:type=\"paragraph\"`,\n# etc.). No text, run/paragraph formatting, style definition, header/\n# footer content, footnote content, section/page setup value, or any\n# other semantically observable property differs between the before\n# and after XML - every w:val / w:id / r:id / dimension stays exactly\n# the same, only the left-to-right order of attributes/xmlns\n# declarations inside tags changes.\n#\n# That kind of change is an artifact of whichever tool re-serialized\n# the .xml parts (e.g. the document being re-saved), not something\n# exposed by the Word COM object model: there is no attribute-order /\n# raw-OOXML-serialization control in the object model, and touching a\n# part (e.g. re-setting styles/headers/footers) to force a re-write\n# would risk introducing real (unwanted) content differences instead\n# of reproducing this purely cosmetic reordering. So there is nothing\n# for this script to change - we just confirm the document is\n# reachable and leave it exactly as authored.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
